# Auto-generated edit script applying numeric updates to leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(3, 8).Value = 0  # H3: was 38885.668
$ws.Cells.Item(3, 10).Value = 0  # J3: was 38885.668
$ws.Cells.Item(3, 12).Value = 0  # L3: was 38885.668
$ws.Cells.Item(3, 14).ClearContents()  # N3: was -39113.668
$ws.Cells.Item(62, 8).Value = 19244836  # H62: was 10424744
$ws.Cells.Item(62, 9).Value = 27791108  # I62: was 13164819
$ws.Cells.Item(62, 10).Value = 15725  # J62: was 12459.8
$ws.Cells.Item(62, 11).Value = 27791108  # K62: was 13164819
$ws.Cells.Item(62, 12).Value = 15725  # L62: was 12459.8
$ws.Cells.Item(62, 13).Value = -27790484  # M62: was -13164195
$ws.Cells.Item(62, 14).Value = -16973  # N62: was -13707.8
$ws.Cells.Item(65, 8).Value = 19244836  # H65: was 10424744
$ws.Cells.Item(65, 9).Value = 27791108  # I65: was 13164819
$ws.Cells.Item(65, 10).Value = 15725  # J65: was 12459.8
$ws.Cells.Item(65, 11).Value = 138955540  # K65: was 65824095
$ws.Cells.Item(65, 12).Value = 78625  # L65: was 62299
$ws.Cells.Item(65, 13).Value = -138952420  # M65: was -65820975
$ws.Cells.Item(65, 14).Value = -84865  # N65: was -68539
$ws.Cells.Item(94, 8).Value = 1835  # H94: was 400
$ws.Cells.Item(94, 9).Value = 1835  # I94: was 400
$ws.Cells.Item(94, 11).Value = 1835  # K94: was 400
$ws.Cells.Item(94, 13).Value = -1384  # M94: was 51
$ws.Cells.Item(100, 8).Value = 14101.25  # H100: was 15915.714
$ws.Cells.Item(100, 9).Value = 14101.25  # I100: was 15915.714
$ws.Cells.Item(100, 11).Value = 14101.25  # K100: was 15915.714
$ws.Cells.Item(100, 13).Value = -13560.25  # M100: was -15374.714
$ws.Cells.Item(102, 8).Value = 0  # H102: was 38885.668
$ws.Cells.Item(102, 10).Value = 0  # J102: was 38885.668
$ws.Cells.Item(102, 12).Value = 0  # L102: was 38885.668
$ws.Cells.Item(102, 14).ClearContents()  # N102: was -45375.668
$ws.Cells.Item(127, 8).Value = 2429.0566  # H127: was 2034.6177
$ws.Cells.Item(127, 9).Value = 495.8  # I127: was 490.45456
$ws.Cells.Item(127, 10).Value = 2878.6511  # J127: was 2773.1304
$ws.Cells.Item(127, 11).Value = 1487.4  # K127: was 1471.36368
$ws.Cells.Item(127, 12).Value = 8635.953300000001  # L127: was 8319.3912
$ws.Cells.Item(127, 13).Value = 3472.6  # M127: was 3488.63632
$ws.Cells.Item(127, 14).Value = -18555.9533  # N127: was -18239.3912

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 144.57143  # H5: was 137.5
$ws.Cells.Item(5, 9).Value = 125.6  # I5: was 133
$ws.Cells.Item(5, 10).Value = 192  # J5: was 140.2
$ws.Cells.Item(5, 11).Value = 125.6  # K5: was 133
$ws.Cells.Item(5, 12).Value = 192  # L5: was 140.2
$ws.Cells.Item(5, 13).Value = -13.59999999999999  # M5: was -21
$ws.Cells.Item(5, 14).Value = -416  # N5: was -364.2
$ws.Cells.Item(44, 8).Value = 43491.5  # H44: was 43658.168
$ws.Cells.Item(44, 10).Value = 43491.5  # J44: was 43658.168
$ws.Cells.Item(44, 12).Value = 43491.5  # L44: was 43658.168
$ws.Cells.Item(44, 14).Value = -44467.5  # N44: was -44634.168
$ws.Cells.Item(55, 8).Value = 23054  # H55: was 23053
$ws.Cells.Item(55, 10).Value = 23054  # J55: was 23053
$ws.Cells.Item(55, 12).Value = 23054  # L55: was 23053
$ws.Cells.Item(55, 14).Value = -23684  # N55: was -23683
$ws.Cells.Item(61, 8).Value = 3789273.2  # H61: was 3473616.8
$ws.Cells.Item(61, 9).Value = 5435703.5  # I61: was 5682839.5
$ws.Cells.Item(61, 10).Value = 2484  # J61: was 1981
$ws.Cells.Item(61, 11).Value = 5435703.5  # K61: was 5682839.5
$ws.Cells.Item(61, 12).Value = 2484  # L61: was 1981
$ws.Cells.Item(61, 13).Value = -5435491.5  # M61: was -5682627.5
$ws.Cells.Item(61, 14).Value = -2908  # N61: was -2405
$ws.Cells.Item(80, 8).Value = 40060  # H80: was 0
$ws.Cells.Item(80, 10).Value = 40060  # J80: was 0
$ws.Cells.Item(80, 12).Value = 40060  # L80: was 0
$ws.Cells.Item(80, 14).Value = -42056  # N80: was None
$ws.Cells.Item(83, 8).Value = 40060  # H83: was 0
$ws.Cells.Item(83, 10).Value = 40060  # J83: was 0
$ws.Cells.Item(83, 12).Value = 120180  # L83: was 0
$ws.Cells.Item(83, 14).Value = -130164  # N83: was None
$ws.Cells.Item(108, 8).Value = 21916  # H108: was 23869.6
$ws.Cells.Item(108, 10).Value = 21916  # J108: was 23869.6
$ws.Cells.Item(108, 12).Value = 21916  # L108: was 23869.6
$ws.Cells.Item(108, 14).Value = -29596  # N108: was -31549.6
$ws.Cells.Item(136, 8).Value = 3789273.2  # H136: was 3473616.8
$ws.Cells.Item(136, 9).Value = 5435703.5  # I136: was 5682839.5
$ws.Cells.Item(136, 10).Value = 2484  # J136: was 1981
$ws.Cells.Item(136, 11).Value = 16307110.5  # K136: was 17048518.5
$ws.Cells.Item(136, 12).Value = 7452  # L136: was 5943
$ws.Cells.Item(136, 13).Value = -16304560.5  # M136: was -17045968.5
$ws.Cells.Item(136, 14).Value = -12552  # N136: was -11043

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 144.57143  # H4: was 137.5
$ws.Cells.Item(4, 9).Value = 125.6  # I4: was 133
$ws.Cells.Item(4, 10).Value = 192  # J4: was 140.2
$ws.Cells.Item(4, 11).Value = 125.6  # K4: was 133
$ws.Cells.Item(4, 12).Value = 192  # L4: was 140.2
$ws.Cells.Item(4, 13).Value = -10.59999999999999  # M4: was -18
$ws.Cells.Item(4, 14).Value = -422  # N4: was -370.2
$ws.Cells.Item(109, 8).Value = 30669  # H109: was 30668.572
$ws.Cells.Item(109, 10).Value = 30669  # J109: was 30668.572
$ws.Cells.Item(109, 12).Value = 30669  # L109: was 30668.572
$ws.Cells.Item(109, 14).Value = -33443  # N109: was -33442.572

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1562  # H16: was 1595.6
$ws.Cells.Item(16, 9).Value = 624.125  # I16: was 448.4
$ws.Cells.Item(16, 10).Value = 2499.875  # J16: was 2169.2
$ws.Cells.Item(16, 11).Value = 624.125  # K16: was 448.4
$ws.Cells.Item(16, 12).Value = 2499.875  # L16: was 2169.2
$ws.Cells.Item(16, 13).Value = -337.125  # M16: was -161.4
$ws.Cells.Item(16, 14).Value = -3073.875  # N16: was -2743.2
$ws.Cells.Item(103, 8).Value = 3334.1428  # H103: was 3934.1428
$ws.Cells.Item(103, 9).Value = 3334.1428  # I103: was 3934.1428
$ws.Cells.Item(103, 11).Value = 3334.1428  # K103: was 3934.1428
$ws.Cells.Item(103, 13).Value = -2162.1428  # M103: was -2762.1428
$ws.Cells.Item(109, 8).Value = 23091.428  # H109: was 23148.143
$ws.Cells.Item(109, 10).Value = 23091.428  # J109: was 23148.143
$ws.Cells.Item(109, 12).Value = 23091.428  # L109: was 23148.143
$ws.Cells.Item(109, 14).Value = -25171.428  # N109: was -25228.143
$ws.Cells.Item(113, 8).Value = 1562  # H113: was 1595.6
$ws.Cells.Item(113, 9).Value = 624.125  # I113: was 448.4
$ws.Cells.Item(113, 10).Value = 2499.875  # J113: was 2169.2
$ws.Cells.Item(113, 11).Value = 624.125  # K113: was 448.4
$ws.Cells.Item(113, 12).Value = 2499.875  # L113: was 2169.2
$ws.Cells.Item(113, 13).Value = 1545.875  # M113: was 1721.6
$ws.Cells.Item(113, 14).Value = -6839.875  # N113: was -6509.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(94, 8).Value = 766.6667  # H94: was 1350
$ws.Cells.Item(94, 9).Value = 766.6667  # I94: was 1350
$ws.Cells.Item(94, 11).Value = 2300.0001  # K94: was 4050
$ws.Cells.Item(94, 13).Value = -1624.0001  # M94: was -3374
$ws.Cells.Item(106, 8).Value = 0  # H106: was 2000
$ws.Cells.Item(106, 10).Value = 0  # J106: was 2000
$ws.Cells.Item(106, 12).Value = 0  # L106: was 6000
$ws.Cells.Item(106, 14).ClearContents()  # N106: was -7892

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 7199.375  # H70: was 11409.571
$ws.Cells.Item(70, 9).Value = 7780.7407  # I70: was 15467.111
$ws.Cells.Item(70, 10).Value = 4060  # J70: was 4106
$ws.Cells.Item(70, 11).Value = 7780.7407  # K70: was 15467.111
$ws.Cells.Item(70, 12).Value = 4060  # L70: was 4106
$ws.Cells.Item(70, 13).Value = -7510.7407  # M70: was -15197.111
$ws.Cells.Item(70, 14).Value = -4600  # N70: was -4646
$ws.Cells.Item(73, 8).Value = 7199.375  # H73: was 11409.571
$ws.Cells.Item(73, 9).Value = 7780.7407  # I73: was 15467.111
$ws.Cells.Item(73, 10).Value = 4060  # J73: was 4106
$ws.Cells.Item(73, 11).Value = 7780.7407  # K73: was 15467.111
$ws.Cells.Item(73, 12).Value = 4060  # L73: was 4106
$ws.Cells.Item(73, 13).Value = -6844.7407  # M73: was -14531.111
$ws.Cells.Item(73, 14).Value = -5932  # N73: was -5978
$ws.Cells.Item(80, 8).Value = 2179.6  # H80: was 2633.3333
$ws.Cells.Item(80, 9).Value = 2133.3333  # I80: was 2700
$ws.Cells.Item(80, 10).Value = 2249  # J80: was 2500
$ws.Cells.Item(80, 11).Value = 2133.3333  # K80: was 2700
$ws.Cells.Item(80, 12).Value = 2249  # L80: was 2500
$ws.Cells.Item(80, 13).Value = -1135.3333  # M80: was -1702
$ws.Cells.Item(80, 14).Value = -4245  # N80: was -4496
$ws.Cells.Item(83, 8).Value = 2179.6  # H83: was 2633.3333
$ws.Cells.Item(83, 9).Value = 2133.3333  # I83: was 2700
$ws.Cells.Item(83, 10).Value = 2249  # J83: was 2500
$ws.Cells.Item(83, 11).Value = 10666.6665  # K83: was 13500
$ws.Cells.Item(83, 12).Value = 11245  # L83: was 12500
$ws.Cells.Item(83, 13).Value = -5674.666499999999  # M83: was -8508
$ws.Cells.Item(83, 14).Value = -21229  # N83: was -22484
$ws.Cells.Item(113, 8).Value = 1241  # H113: was 1623.4445
$ws.Cells.Item(113, 9).Value = 1183.3  # I113: was 1544.4286
$ws.Cells.Item(113, 10).Value = 1433.3334  # J113: was 1900
$ws.Cells.Item(113, 11).Value = 1183.3  # K113: was 1544.4286
$ws.Cells.Item(113, 12).Value = 1433.3334  # L113: was 1900
$ws.Cells.Item(113, 13).Value = 986.7  # M113: was 625.5714
$ws.Cells.Item(113, 14).Value = -5773.3334  # N113: was -6240

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 2522  # H100: was 2618.1
$ws.Cells.Item(100, 9).Value = 1841  # I100: was 1996.125
$ws.Cells.Item(100, 11).Value = 1841  # K100: was 1996.125
$ws.Cells.Item(100, 13).Value = -1300  # M100: was -1455.125
$ws.Cells.Item(132, 8).Value = 21170566  # H132: was 24323514
$ws.Cells.Item(132, 9).Value = 31748080  # I132: was 40818732
$ws.Cells.Item(132, 10).Value = 15539.723  # J132: was 14769.211
$ws.Cells.Item(132, 11).Value = 95244240  # K132: was 122456196
$ws.Cells.Item(132, 12).Value = 46619.169  # L132: was 44307.633
$ws.Cells.Item(132, 13).Value = -95241710  # M132: was -122453666
$ws.Cells.Item(132, 14).Value = -51679.169  # N132: was -49367.633

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 4293.8887  # H62: was 4380.25
$ws.Cells.Item(62, 10).Value = 4801.5  # J62: was 6000
$ws.Cells.Item(62, 12).Value = 4801.5  # L62: was 6000
$ws.Cells.Item(62, 14).Value = -6049.5  # N62: was -7248
$ws.Cells.Item(65, 8).Value = 4293.8887  # H65: was 4380.25
$ws.Cells.Item(65, 10).Value = 4801.5  # J65: was 6000
$ws.Cells.Item(65, 12).Value = 24007.5  # L65: was 30000
$ws.Cells.Item(65, 14).Value = -30247.5  # N65: was -36240
$ws.Cells.Item(109, 8).Value = 15340.333  # H109: was 20341.8
$ws.Cells.Item(109, 10).Value = 15340.333  # J109: was 20341.8
$ws.Cells.Item(109, 12).Value = 15340.333  # L109: was 20341.8
$ws.Cells.Item(109, 14).Value = -18114.333  # N109: was -23115.8
